# Add a new "health-costs" worksheet between "parameters" and "jurisdiction",
# matching the data added in the target workbook.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("parameters")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "health-costs"

# Populate row-by-row, bottom-to-top, right-to-left, matching the order the
# shared strings were originally authored in.

$newSheet.Range("A5").Value = "chronic"
$newSheet.Range("B5").Value = 0.219
$newSheet.Range("C5").Value = 1
$newSheet.Range("D5").Value = 0.05
$newSheet.Range("E5").Value = 279113
$newSheet.Range("F5").Value = 0

$newSheet.Range("A4").Value = "critical"
$newSheet.Range("B4").Value = 0.655
$newSheet.Range("C4").Value = 0.036
$newSheet.Range("D4").Value = 0.03
$newSheet.Range("E4").Value = 279113
$newSheet.Range("F4").Value = 41510

$newSheet.Range("A3").Value = "severe"
$newSheet.Range("B3").Value = 0.133
$newSheet.Range("C3").Value = 0.023
$newSheet.Range("D3").Value = 0.09
$newSheet.Range("E3").Value = 279113
$newSheet.Range("F3").Value = 11267

$newSheet.Range("A2").Value = "mild"
$newSheet.Range("B2").Value = 0.006
$newSheet.Range("C2").Value = 0.013
$newSheet.Range("D2").Value = 0.54
$newSheet.Range("E2").Value = 279113
$newSheet.Range("F2").Value = 0

$newSheet.Range("F1").Value = "hospital_cost"
$newSheet.Range("E1").Value = "VSLY"
$newSheet.Range("D1").Value = "disease_state_prevalence"
$newSheet.Range("C1").Value = "disease_duration"
$newSheet.Range("B1").Value = "DALY_weight"

# Column widths to match the author's table formatting (bestFit autosized
# columns). The saved/stored column width ends up 5/6 of a character wider
# than the ColumnWidth value we assign here, so back that constant out to
# land on the target stored widths (7, 12, 15.16..., 22.66..., 7.16..., 12).
$widthOffset = 5 / 6
$newSheet.Columns.Item(1).ColumnWidth = 7 - $widthOffset
$newSheet.Columns.Item(2).ColumnWidth = 12 - $widthOffset
$newSheet.Columns.Item(3).ColumnWidth = 15.1640625 - $widthOffset
$newSheet.Columns.Item(4).ColumnWidth = 22.6640625 - $widthOffset
$newSheet.Columns.Item(5).ColumnWidth = 7.1640625 - $widthOffset
$newSheet.Columns.Item(6).ColumnWidth = 12 - $widthOffset

# Select the same cell the author left active on the new tab.
$newSheet.Range("F28").Select()

# Make the new tab the active one (mirrors the author leaving health-costs
# selected, and bumping the previously-active "jurisdiction" tab along).
$newSheet.Activate()
